# NIT-9007803599.xlsx — "Elimina antiguos EC y agrega nuevos y modifica Antigua BD"
#
# Adds a new "Estado de Cuenta" (EC) employee row to the arrears table
# (IVAN DARIO CASTILLO GARAY) and refreshes the summary totals on the
# existing worksheet accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets("Hoja1")

# Insert a new row right after the existing employee record (row 16) so the
# footer rows (signature block) shift down by one, same as the target file.
$ws.Rows("17:17").Insert()

# Clone the formatting of the existing data row into the freshly inserted
# one so the new record matches the table's look (borders, fills, number
# formats, etc.).
$ws.Range("B16:J16").Copy()
$ws.Range("B17:J17").PasteSpecial(-4122)

# New employee arrears record.
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "18859019"
$ws.Range("D17").Value = "IVAN DARIO CASTILLO GARAY"
$ws.Range("E17").Value = "2509"
$ws.Range("F17").Value = 22776
$ws.Range("G17").Value = 1423500

# Refresh the header summary: total "VALOR MORA" (25396 + 22776), and the
# worker/period counters (now 2 each instead of 1).
$ws.Range("E11").Value = 48172
$ws.Range("C13").Value = 2
$ws.Range("F13").Value = 2
